$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44306
$ws.Range("M2").Value2 = 80

# Row 3
$ws.Range("D3").Value2 = 44316

# Row 4
$ws.Range("D4").Value2 = 44309
$ws.Range("Q4").Value = "$/caja 14 kilos granel"
$ws.Range("S4").Value2 = 821
$ws.Range("T4").Value2 = 14

# Row 5
$ws.Range("D5").Value2 = 44313

# Row 6
$ws.Range("D6").Value2 = 44330

# Row 7
$ws.Range("D7").Value2 = 44322
$ws.Range("M7").Value2 = 60

# Row 8
$ws.Range("D8").Value2 = 44302

# Row 9
$ws.Range("D9").Value2 = 44323
$ws.Range("M9").Value2 = 80

# Row 10
$ws.Range("D10").Value2 = 44327
$ws.Range("M10").Value2 = 60
$ws.Range("Q10").Value = "$/caja 10 kilos empedrada"
$ws.Range("S10").Value2 = 11500
$ws.Range("T10").Value2 = 1
